$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header in C1: "Refcode" -> "Emails Refcode"
$ws.Range("C1").Value = "Emails Refcode"

# Move the active selection to D1
$ws.Range("D1").Select()
